$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number + date range)
$ws.Range("A8").Value = "Volume 29   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/7/2022  Through  11/13/2022"

$ws.Range("C14").NumberFormat = "General"
$ws.Range("C14").Value = "'0"

$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 28
$ws.Range("K15").Value = 12
$ws.Range("L15").Value = 3.703703703703
$ws.Range("M15").Value = 55.555555555555
$ws.Range("N15").Value = -61.643835616438

$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = 15.384615384615
$ws.Range("F16").Value = 39
$ws.Range("G16").Value = 48
$ws.Range("H16").Value = -18.75
$ws.Range("I16").Value = 489
$ws.Range("J16").Value = 370
$ws.Range("K16").Value = 32.162162162162
$ws.Range("L16").Value = 26.356589147286
$ws.Range("M16").Value = 23.173803526448
$ws.Range("N16").Value = -70.718562874251

$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -13.333333333333
$ws.Range("F17").Value = 61
$ws.Range("G17").Value = 62
$ws.Range("H17").Value = -1.612903225806
$ws.Range("I17").Value = 739
$ws.Range("J17").Value = 611
$ws.Range("K17").Value = 20.949263502455
$ws.Range("L17").Value = 15.830721003134
$ws.Range("M17").Value = 93.455497382199
$ws.Range("N17").Value = -24.97461928934

$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 33
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 323
$ws.Range("J18").Value = 162
$ws.Range("K18").Value = 99.382716049382
$ws.Range("L18").Value = 30.76923076923
$ws.Range("M18").Value = 83.522727272727
$ws.Range("N18").Value = -72.650296359017

$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -38.888888888888
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 70
$ws.Range("H19").Value = -35.714285714285
$ws.Range("I19").Value = 642
$ws.Range("J19").Value = 638
$ws.Range("K19").Value = 0.626959247648
$ws.Range("L19").Value = 2.229299363057
$ws.Range("M19").Value = 69.841269841269
$ws.Range("N19").Value = -4.606240713224

$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 273
$ws.Range("J20").Value = 159
$ws.Range("K20").Value = 71.698113207547
$ws.Range("L20").Value = 79.605263157894
$ws.Range("M20").Value = 155.140186915888
$ws.Range("N20").Value = -53.173241852487

$ws.Range("C21").Value = 52
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = -5.454545454545
$ws.Range("F21").Value = 203
$ws.Range("G21").Value = 219
$ws.Range("H21").Value = -7.305936073059
$ws.Range("I21").Value = 2508
$ws.Range("J21").Value = 1980
$ws.Range("K21").Value = 26.666666666666
$ws.Range("L21").Value = 19.942611190817
$ws.Range("M21").Value = 70.380434782608
$ws.Range("N21").Value = -52.027543993879

$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 78
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 30
$ws.Range("L22").Value = 13.043478260869
$ws.Range("M22").Value = 47.169811320754

$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 166.666666666667
$ws.Range("F23").Value = 29
$ws.Range("G23").Value = 26
$ws.Range("H23").Value = 11.538461538461
$ws.Range("I23").Value = 345
$ws.Range("J23").Value = 251
$ws.Range("K23").Value = 37.450199203187
$ws.Range("L23").Value = 26.373626373626
$ws.Range("M23").Value = 31.679389312977

$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -21.212121212121
$ws.Range("F24").Value = 164
$ws.Range("G24").Value = 121
$ws.Range("H24").Value = 35.537190082644
$ws.Range("I24").Value = 1677
$ws.Range("J24").Value = 1180
$ws.Range("K24").Value = 42.118644067796
$ws.Range("L24").Value = 18.599717114568
$ws.Range("M24").Value = 36.897959183673

$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -42.307692307692
$ws.Range("F25").Value = 75
$ws.Range("G25").Value = 83
$ws.Range("H25").Value = -9.638554216867
$ws.Range("I25").Value = 903
$ws.Range("J25").Value = 806
$ws.Range("K25").Value = 12.034739454094
$ws.Range("L25").Value = -2.588996763754
$ws.Range("M25").Value = -10.94674556213

$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 1
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E26").Value = "***.*"
$ws.Range("I26").Value = 46
$ws.Range("K26").Value = -25.806451612903
$ws.Range("L26").Value = 6.976744186046

$ws.Range("C27").NumberFormat = "General"
$ws.Range("C27").Value = "'0"
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 76
$ws.Range("K27").Value = -2.631578947368

$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -55.555555555555
$ws.Range("I28").Value = 56
$ws.Range("J28").Value = 64
$ws.Range("K28").Value = -12.5
$ws.Range("L28").Value = -15.151515151515
$ws.Range("M28").Value = 1.818181818181
$ws.Range("N28").Value = -71.573604060913

$ws.Range("C29").Value = 1
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = -62.5
$ws.Range("I29").Value = 46
$ws.Range("J29").Value = 56
$ws.Range("K29").Value = -17.857142857142
$ws.Range("L29").Value = -13.207547169811
$ws.Range("M29").Value = -2.127659574468
$ws.Range("N29").Value = -74.011299435028

$ws.Range("C30").NumberFormat = "General"
$ws.Range("C30").Value = "'0"
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -100
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("G30").Value = 1
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = -20
